$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its original text data type,
# since Excel would otherwise auto-convert numeric-looking strings
# (e.g. "1.007") into actual numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.471.89"
$ws.Range("E2").Value = "  -3.30%  "
$ws.Range("D3").Value = "1.754.76"
$ws.Range("E3").Value = "  -2.77%  "
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  +0.56%  "
$ws.Range("D5").Value = "324.54"
$ws.Range("E5").Value = "  -0.94%  "
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("D7").Value = "0.4420"
$ws.Range("E7").Value = "  -0.46%  "
$ws.Range("D8").Value = "0.3687"
$ws.Range("E8").Value = "  -2.42%  "
$ws.Range("D9").Value = "44.86"
$ws.Range("D10").Value = "0.07649"
$ws.Range("E10").Value = "  +1.90%  "
$ws.Range("D11").Value = "1.110"
$ws.Range("E11").Value = "  -3.51%  "
$ws.Range("D12").Value = "1.006"
$ws.Range("E12").Value = "  +0.60%  "
$ws.Range("D13").Value = "21.55"
$ws.Range("E13").Value = "  -4.61%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "6.140"
$ws.Range("E14").Value = "  -2.51%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "7.420"
$ws.Range("E15").Value = "  -3.01%  "
$ws.Range("D16").Value = "1.758.17"
$ws.Range("E16").Value = "  -2.28%  "
$ws.Range("D17").Value = "90.01"
$ws.Range("E17").Value = "  +11.50%  "
$ws.Range("D18").Value = "0.00001068"
$ws.Range("E18").Value = "  -2.17%  "
$ws.Range("D19").Value = "0.06257"
$ws.Range("E19").Value = "  -8.07%  "
$ws.Range("D20").Value = "1.003"
$ws.Range("E20").Value = "  +0.35%  "
$ws.Range("D21").Value = "17.36"
$ws.Range("E21").Value = "  -0.73%  "
$ws.Range("D22").Value = "6.158"
$ws.Range("E22").Value = "  -2.50%  "
$ws.Range("D23").Value = "0.5318"
$ws.Range("E23").Value = "  -2.56%  "
$ws.Range("D24").Value = "27.528.11"
$ws.Range("E24").Value = "  -3.00%  "
$ws.Range("D25").Value = "11.51"
$ws.Range("E25").Value = "  -2.48%  "
$ws.Range("D26").Value = "2.310"
$ws.Range("E26").Value = "  -3.83%  "
$ws.Range("D27").Value = "20.50"
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("D28").Value = "152.68"
$ws.Range("E28").Value = "  -0.37%  "
$ws.Range("D29").Value = "2.281"
$ws.Range("E29").Value = "  -3.51%  "
$ws.Range("D30").Value = "1.961.14"
$ws.Range("E30").Value = "  -2.30%  "
$ws.Range("D31").Value = "127.70"
$ws.Range("E31").Value = "  -3.55%  "
$ws.Range("D32").Value = "1.171"
$ws.Range("E32").Value = "  -6.79%  "
$ws.Range("D33").Value = "5.692"
$ws.Range("E33").Value = "  -2.20%  "
$ws.Range("D34").Value = "0.09185"
$ws.Range("E34").Value = "  -1.40%  "
$ws.Range("D35").Value = "3.622"
$ws.Range("E35").Value = "  -9.51%  "
$ws.Range("D36").Value = "12.54"
$ws.Range("E36").Value = "  +3.43%  "
$ws.Range("D37").Value = "0.02305"
$ws.Range("E37").Value = "  -0.78%  "
$ws.Range("D38").Value = "0.2152"
$ws.Range("E38").Value = "  -5.99%  "
$ws.Range("D39").Value = "0.06085"
$ws.Range("E39").Value = "  -4.37%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "0.6413"
$ws.Range("E40").Value = "  -2.71%  "
$ws.Range("B41").Value = "InternetComputer(DFINITY)"
$ws.Range("C41").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D41").Value = "5.032"
$ws.Range("E41").Value = "  -2.58%  "
$ws.Range("D42").Value = "1.171"
$ws.Range("E42").Value = "  -3.09%  "
$ws.Range("B43").Value = "Frax"
$ws.Range("C43").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D43").Value = "1.003"
$ws.Range("E43").Value = "  +0.32%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "7.921"
$ws.Range("E44").Value = "  -2.87%  "
$ws.Range("D45").Value = "1.394"
$ws.Range("E45").Value = "  -4.05%  "
$ws.Range("E46").Value = "  -1.38%  "
$ws.Range("D47").Value = "0.5946"
$ws.Range("E47").Value = "  -2.10%  "
$ws.Range("D48").Value = "3.719"
$ws.Range("E48").Value = "  -2.32%  "
$ws.Range("D49").Value = "125.95"
$ws.Range("E49").Value = "  -2.02%  "
$ws.Range("D50").Value = "1.983"
$ws.Range("E50").Value = "  -2.52%  "
$ws.Range("D51").Value = "0.06886"
$ws.Range("E51").Value = "  -2.84%  "
